$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (index 3), rows 2-25
$CValues = @(0.005662360753353468, 0.005464652714179508, 0.005341410750705933, 0.005290722662971703, 0.005282277766172072, 0.005340729042476511, 0.005594573897214516, 0.006077769785910903, 0.006424007598560877, 0.006579646417637974, 0.006638316137816958, 0.006625692435761721, 0.006584478562636775, 0.006559199093750578, 0.00641379876855197, 0.006324123084787914, 0.006272368254670369, 0.006254814741289749, 0.006333687396491428, 0.006596591322274747, 0.006766856577542057, 0.00667612497350234, 0.006329363991838477, 0.00594860138151887)
for ($i = 0; $i -lt $CValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $CValues[$i]
}

# Column D (index 4), rows 2-25
$DValues = @(0.1657542215531862, 0.1656575482198335, 0.1656594146240096, 0.1656755819618994, 0.1656791974117198, 0.1656595702656389, 0.1657081878803552, 0.16628899260526, 0.1670115016796458, 0.1674044103069079, 0.1675624270675939, 0.167527984882625, 0.1674172254760649, 0.1673505840866483, 0.1669871162197865, 0.1667805881663753, 0.1666678451968977, 0.1666307110545659, 0.1668019476872118, 0.1674495077245268, 0.1679265272181212, 0.1676670120150376, 0.1667922723812225, 0.166079915444044)
for ($i = 0; $i -lt $DValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $DValues[$i]
}

# Column E (index 5), rows 2-25
$EValues = @(0.1560617972850373, 0.1538024989725777, 0.15250064878731, 0.1519915571079054, 0.1519083154531558, 0.1524936963151653, 0.1552650454768632, 0.1613794864529368, 0.166290372562834, 0.1686163153965623, 0.1695103775833076, 0.1693172338602551, 0.1686896039987644, 0.1683068934385759, 0.1661402217915224, 0.1648346298184933, 0.1640923398405079, 0.1638424975403652, 0.164972716469471, 0.168873593295487, 0.1715004631242394, 0.1700913507559534, 0.1649102616081919, 0.1596521609232973)
for ($i = 0; $i -lt $EValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $EValues[$i]
}

# Column F (index 6), rows 2-25
$FValues = @(1.671014047417117, 1.618929075160096, 1.587975245859468, 1.575618501119649, 1.573582179559622, 1.587807558372148, 1.65284149624074, 1.788569721226452, 1.893370861341026, 1.942170287001204, 1.960812432054894, 1.956790251411718, 1.943700717094629, 1.935704241955278, 1.890204433495285, 1.862580661376484, 1.846798008239034, 1.841472412555532, 1.865510300207319, 1.947541001679895, 2.002102616583727, 1.97289477741046, 1.864185503376831, 1.750965426991968)
for ($i = 0; $i -lt $FValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $FValues[$i]
}

# Column G (index 7), rows 2-25
$GValues = @(0.002420309550770803, 0.002425096453804731, 0.002428193274214087, 0.00242949502869159, 0.002429713590488196, 0.002428210668750603, 0.002421927437248006, 0.002410850731554764, 0.002403462898412298, 0.002400263057641723, 0.00239907436153706, 0.002399329347049469, 0.002400164802219066, 0.002400679536398229, 0.002403675242834219, 0.002405554138076697, 0.002406649982815708, 0.002407023623589373, 0.002405352559034651, 0.002399918785288443, 0.002396501587220851, 0.002398313182588902, 0.002405443644082295, 0.002413714909989206)
for ($i = 0; $i -lt $GValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $GValues[$i]
}

# Column I (index 9), rows 2-25
$IValues = @(0.9185472614003345, 0.8822415908493184, 0.860529189211249, 0.8518262437820994, 0.8503898667386096, 0.860411231891959, 0.9059083892653348, 0.9997616778198903, 1.071598801724491, 1.104918649237192, 1.11762906363505, 1.114887503869923, 1.105962476160215, 1.100507761758436, 1.069434230842447, 1.050536287986034, 1.039726957704872, 1.036077431828176, 1.052541763386913, 1.108581445382029, 1.145748539673903, 1.125861919307937, 1.051634915608886, 0.9738698291636183)
for ($i = 0; $i -lt $IValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 9).Value = $IValues[$i]
}

# Column J (index 10), rows 2-25
$JValues = @(0.1932238490453386, 0.1878352636291609, 0.1846464827536138, 0.1833770359118958, 0.1831680539588589, 0.1846292412152337, 0.1913409155047248, 0.2054599256648686, 0.2164273938263221, 0.22154819456604, 0.2235064070283386, 0.2230838201640353, 0.2217089147212761, 0.2208692341137919, 0.2160953989148879, 0.2132006262025925, 0.2115480116453909, 0.2109905869807989, 0.2135074971418334, 0.2221122388767611, 0.2278471959355812, 0.2247761122314671, 0.2133687245810449, 0.2015367754054154)
for ($i = 0; $i -lt $JValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 10).Value = $JValues[$i]
}

# Column N (index 14), rows 2-25
$NValues = @(2.873581687696344, 2.562605684679454, 2.371325805375761, 2.293303068607429, 2.280343261403573, 2.370273851395496, 2.766433886209654, 3.540180268007646, 4.10623028343673, 4.363110593465422, 4.460285735714251, 4.439361943450422, 4.371107314139522, 4.329286057409945, 4.089429168003562, 3.94211849063862, 3.857331695637754, 3.828614786364199, 3.957806003280837, 4.391158149571254, 4.673791817957863, 4.523002190001307, 3.950713976768498, 3.331249627311138)
for ($i = 0; $i -lt $NValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 14).Value = $NValues[$i]
}

# Column O (index 15), rows 2-25
$OValues = @(4.250172472561133, 4.09312108477036, 3.999359741489343, 3.961819388656693, 3.955626111343065, 3.998850757741138, 4.195465394021312, 4.602361170041831, 4.914571738137852, 5.059543080872459, 5.114867506673932, 5.102933359247345, 5.064086078761079, 5.040346698035592, 4.905157075009072, 4.822979363760396, 4.775989925992917, 4.760127583089343, 4.831698627653395, 5.075484852666364, 5.237303306422689, 5.150708771068309, 4.827755852455084, 4.489974540659034)
for ($i = 0; $i -lt $OValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 15).Value = $OValues[$i]
}
